# Regenerate the ObjTables workbook: bump the embedded generation
# timestamp in each sheet's "!!...ObjTables ... date='...'" marker
# string, and populate the previously-blank "Verbose name" column (E)
# on the Schema sheet for every Attribute row.

$wb = $excel.ActiveWorkbook

$newDate = "2020-05-29 00:18:59"

$toc          = $wb.Worksheets.Item("!!_Table of contents")
$schema       = $wb.Worksheets.Item("!!_Schema")
$genes        = $wb.Worksheets.Item("!!Genes")
$transcripts  = $wb.Worksheets.Item("!!Transcripts")

# --- Table of contents sheet: A1 and A2 hold locked header cells -------
$toc.Range("A1").Locked = $false
$toc.Range("A2").Locked = $false
$toc.Range("A1").Value = "!!!ObjTables objTablesVersion='1.0.0' date='$newDate'"
$toc.Range("A2").Value = "!!ObjTables type='TableOfContents' tableFormat='row' description='Table of contents' date='$newDate' objTablesVersion='1.0.0'"
$toc.Range("A1").Locked = $true
$toc.Range("A2").Locked = $true

# --- Schema sheet: A1 locked header cell --------------------------------
$schema.Range("A1").Locked = $false
$schema.Range("A1").Value = "!!ObjTables type='Schema' tableFormat='row' description='Table/model and column/attribute definitions' date='$newDate' objTablesVersion='1.0.0'"
$schema.Range("A1").Locked = $true

# Fill in the "Verbose name" column (E) for every Attribute row; these
# cells are unlocked data cells so no protection toggling is needed.
$schema.Range("E4").Value = "Id"
$schema.Range("E5").Value = "Location"
$schema.Range("E6").Value = "Symbol"
$schema.Range("E8").Value = "Gene"
$schema.Range("E9").Value = "Id"
$schema.Range("E10").Value = "Location"
$schema.Range("E12").Value = "Chromosome"
$schema.Range("E13").Value = "5'"
$schema.Range("E14").Value = "3'"

# --- Genes sheet: A1 locked header cell ---------------------------------
$genes.Range("A1").Locked = $false
$genes.Range("A1").Value = "!!ObjTables type='Data' tableFormat='row' class='Gene' name='Genes' date='$newDate' objTablesVersion='1.0.0'"
$genes.Range("A1").Locked = $true

# --- Transcripts sheet: A1 locked header cell ---------------------------
$transcripts.Range("A1").Locked = $false
$transcripts.Range("A1").Value = "!!ObjTables type='Data' tableFormat='row' class='Transcript' name='Transcripts' date='$newDate' objTablesVersion='1.0.0'"
$transcripts.Range("A1").Locked = $true
